# Push first 11 samples with replicates from 20190112 20 degrees light samples Moorea
# (two new CRM-accuracy rows, 27 and 28, each with a Date/Batch-value/Note entry)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 -------------------------------------------------------------
# Copy A26's format (date-formatted style) down into A27 before writing the
# value so the new cell reuses the existing date style instead of minting a
# new one.
$ws.Range("A26").Copy($ws.Range("A27"))
$ws.Range("A27").Value = 43503
$ws.Range("B27").Value = 2189.90222552467
$ws.Range("F27").Value = "CRM bottle already opened "

# --- Row 28 -------------------------------------------------------------
$ws.Range("A26").Copy($ws.Range("A28"))
$ws.Range("A28").Value = 43503
$ws.Range("B28").Value = 2204.7327702513298
$ws.Range("F28").Value = "New CRM bottle"

# --- New column widths for the note columns (F, G) ----------------------
$ws.Columns.Item(6).ColumnWidth = 27.333333333333332
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334

# --- View state: scroll position + active selection ---------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("F29").Select()
